$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 13412
$ws.Range("E2").Value = 404
$ws.Range("F2").Value = 492
$ws.Range("G2").Value = 2984
$ws.Range("H2").Value = 2905
$ws.Range("I2").Value = 2907
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = 25230
$ws.Range("L2").Value = 10492
$ws.Range("M2").Value = 14738
$ws.Range("N2").Value = 13694
$ws.Range("O2").Value = 1044
$ws.Range("P2").Value = 150
$ws.Range("Q2").Value = 746
$ws.Range("R2").Value = -644
$ws.Range("S2").Value = 1290
$ws.Range("T2").Value = 1286
$ws.Range("U2").Value = -540
$ws.Range("V2").Value = 6400
$ws.Range("W2").Value = 3.01
$ws.Range("X2").Value = 21.66
$ws.Range("Y2").Value = 23.02
$ws.Range("Z2").Value = 13.37
$ws.Range("AA2").Value = 71.19
$ws.Range("AB2").Value = 8986.469999999999
$ws.Range("AC2").Value = 14535
$ws.Range("AD2").Value = 2.62
$ws.Range("AE2").Value = 68481
$ws.Range("AF2").Value = 0.5600000000000001
$ws.Range("AG2").Value = 262
$ws.Range("AH2").Value = 0.6899999999999999
$ws.Range("AI2").Value = 1.81
$ws.Range("AJ2").Value = 20000000

# Row 3
$ws.Range("D3").Value = 14732
$ws.Range("E3").Value = 444
$ws.Range("F3").Value = 431
$ws.Range("G3").Value = 291
$ws.Range("H3").Value = 183
$ws.Range("I3").Value = 194
$ws.Range("J3").Value = -10
$ws.Range("K3").Value = 24258
$ws.Range("L3").Value = 9698
$ws.Range("M3").Value = 14560
$ws.Range("N3").Value = 13560
$ws.Range("O3").Value = 1000
$ws.Range("P3").Value = 180
$ws.Range("Q3").Value = 2007
$ws.Range("R3").Value = -1554
$ws.Range("S3").Value = -668
$ws.Range("T3").Value = 1778
$ws.Range("U3").Value = 230
$ws.Range("V3").Value = 5772
$ws.Range("W3").Value = 3.01
$ws.Range("X3").Value = 1.24
$ws.Range("Y3").Value = 1.42
$ws.Range("Z3").Value = 0.74
$ws.Range("AA3").Value = 66.59999999999999
$ws.Range("AB3").Value = 7519.24
$ws.Range("AC3").Value = 968
$ws.Range("AD3").Value = 35.43
$ws.Range("AE3").Value = 67816
$ws.Range("AF3").Value = 0.51
$ws.Range("AG3").Value = 315
$ws.Range("AH3").Value = 0.92
$ws.Range("AI3").Value = 32.52
$ws.Range("AJ3").Value = 20000000

# Row 4
$ws.Range("D4").Value = 14174
$ws.Range("E4").Value = 370
$ws.Range("F4").Value = 384
$ws.Range("G4").Value = 377
$ws.Range("H4").Value = 403
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 24531
$ws.Range("L4").Value = 9708
$ws.Range("M4").Value = 14823
$ws.Range("N4").Value = 13850
$ws.Range("O4").Value = 974
$ws.Range("P4").Value = 180
$ws.Range("Q4").Value = 1192
$ws.Range("R4").Value = -1570
$ws.Range("S4").Value = -114
$ws.Range("T4").Value = 1576
$ws.Range("U4").Value = -384
$ws.Range("V4").Value = 5784
$ws.Range("W4").Value = 2.61
$ws.Range("X4").Value = 2.84
$ws.Range("Y4").Value = 2.92
$ws.Range("Z4").Value = 1.65
$ws.Range("AA4").Value = 65.48999999999999
$ws.Range("AB4").Value = 7718.04
$ws.Range("AC4").Value = 1999
$ws.Range("AD4").Value = 18.09
$ws.Range("AE4").Value = 69262
$ws.Range("AF4").Value = 0.52
$ws.Range("AG4").Value = 315
$ws.Range("AH4").Value = 0.87
$ws.Range("AI4").Value = 15.75
$ws.Range("AJ4").Value = 20000000

# Row 5
$ws.Range("D5").Value = 15505
$ws.Range("E5").Value = 388
$ws.Range("F5").Value = 388
$ws.Range("G5").Value = 431
$ws.Range("H5").Value = 225
$ws.Range("I5").Value = 237
$ws.Range("J5").Value = -12
$ws.Range("K5").Value = 24006
$ws.Range("L5").Value = 9423
$ws.Range("M5").Value = 14583
$ws.Range("N5").Value = 13656
$ws.Range("O5").Value = 926
$ws.Range("P5").Value = 180
$ws.Range("Q5").Value = 1140
$ws.Range("R5").Value = -936
$ws.Range("S5").Value = -180
$ws.Range("T5").Value = 974
$ws.Range("U5").Value = 166
$ws.Range("V5").Value = 5261
$ws.Range("W5").Value = 2.5
$ws.Range("X5").Value = 1.45
$ws.Range("Y5").Value = 1.72
$ws.Range("Z5").Value = 0.93
$ws.Range("AA5").Value = 64.61
$ws.Range("AB5").Value = 7827.65
$ws.Range("AC5").Value = 1184
$ws.Range("AD5").Value = 27.69
$ws.Range("AE5").Value = 68296
$ws.Range("AF5").Value = 0.48
$ws.Range("AG5").Value = 315
$ws.Range("AH5").Value = 0.96
$ws.Range("AI5").Value = 26.59
$ws.Range("AJ5").Value = 20000000

# Row 6
$ws.Range("D6").Value = 15772
$ws.Range("E6").Value = 380
$ws.Range("F6").Value = 380
$ws.Range("G6").Value = 79
$ws.Range("H6").Value = 19
$ws.Range("I6").Value = 34
$ws.Range("K6").Value = 23875
$ws.Range("L6").Value = 9304
$ws.Range("M6").Value = 14571
$ws.Range("N6").Value = 13640
$ws.Range("P6").Value = 200
$ws.Range("Q6").Value = 1240
$ws.Range("R6").Value = -526
$ws.Range("S6").Value = -385
$ws.Range("T6").Value = 805
$ws.Range("U6").Value = 435
$ws.Range("V6").Value = 4909
$ws.Range("W6").Value = 2.41
$ws.Range("X6").Value = 0.12
$ws.Range("Y6").Value = 0.25
$ws.Range("Z6").Value = 0.08
$ws.Range("AA6").Value = 63.85
$ws.Range("AB6").Value = 7026.93
$ws.Range("AC6").Value = 169
$ws.Range("AD6").Value = 137.43
$ws.Range("AE6").Value = 68219
$ws.Range("AF6").Value = 0.34
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 1.51
$ws.Range("AI6").Value = 206.82
$ws.Range("AJ6").Value = 20000000

# Clear rows 7-9 (D:AJ), only keep columns A-C
$ws.Range("D7:AJ9").ClearContents()
